$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '30.410.17'
$ws.Cells.Item(2, 5).Value = '  +0.38%  '
$ws.Cells.Item(3, 4).Value = '1.875.02'
$ws.Cells.Item(3, 5).Value = '  -0.23%  '
Set-TextValue $ws.Cells.Item(4, 4) '0.9999'
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
Set-TextValue $ws.Cells.Item(5, 4) '238.06'
$ws.Cells.Item(5, 5).Value = '  +0.72%  '
$ws.Cells.Item(6, 5).Value = '  +0.11%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.4826'
$ws.Cells.Item(7, 5).Value = '  -0.24%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.2823'
$ws.Cells.Item(8, 5).Value = '  -1.93%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.06548'
$ws.Cells.Item(9, 5).Value = '  -0.68%  '
$ws.Cells.Item(10, 4).Value = '1.885.07'
$ws.Cells.Item(10, 5).Value = '  +0.27%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.07448'
$ws.Cells.Item(11, 5).Value = '  +2.22%  '
Set-TextValue $ws.Cells.Item(12, 4) '16.49'
$ws.Cells.Item(12, 5).Value = '  -1.49%  '
Set-TextValue $ws.Cells.Item(13, 4) '5.086'
$ws.Cells.Item(13, 5).Value = '  -2.00%  '
Set-TextValue $ws.Cells.Item(14, 4) '88.01'
$ws.Cells.Item(14, 5).Value = '  +0.88%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.6574'
$ws.Cells.Item(15, 5).Value = '  +0.21%  '
$ws.Cells.Item(16, 4).Value = '30.377.14'
$ws.Cells.Item(16, 5).Value = '  +0.43%  '
Set-TextValue $ws.Cells.Item(17, 4) '13.34'
$ws.Cells.Item(17, 5).Value = '  -0.12%  '
Set-TextValue $ws.Cells.Item(18, 4) '0.9993'
$ws.Cells.Item(18, 5).Value = '  -0.01%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.000007634'
$ws.Cells.Item(19, 5).Value = '  -0.96%  '
$ws.Cells.Item(20, 4).Value = '2.120.11'
$ws.Cells.Item(20, 5).Value = '  -0.16%  '
Set-TextValue $ws.Cells.Item(21, 4) '5.296'
$ws.Cells.Item(21, 5).Value = '  -0.32%  '
Set-TextValue $ws.Cells.Item(22, 4) '0.9996'
$ws.Cells.Item(22, 5).Value = '  +0.06%  '
Set-TextValue $ws.Cells.Item(23, 4) '221.40'
$ws.Cells.Item(23, 5).Value = '  +13.42%  '
Set-TextValue $ws.Cells.Item(24, 4) '6.193'
$ws.Cells.Item(24, 5).Value = '  +1.15%  '
Set-TextValue $ws.Cells.Item(25, 4) '9.249'
$ws.Cells.Item(25, 5).Value = '  -0.50%  '
Set-TextValue $ws.Cells.Item(26, 4) '165.00'
$ws.Cells.Item(26, 5).Value = '  +3.87%  '
Set-TextValue $ws.Cells.Item(27, 4) '18.55'
$ws.Cells.Item(27, 5).Value = '  +2.61%  '
Set-TextValue $ws.Cells.Item(28, 4) '1.980'
$ws.Cells.Item(28, 5).Value = '  +3.33%  '
Set-TextValue $ws.Cells.Item(29, 4) '1.454'
$ws.Cells.Item(29, 5).Value = '  +0.74%  '
Set-TextValue $ws.Cells.Item(30, 4) '0.09402'
$ws.Cells.Item(30, 5).Value = '  +2.90%  '
Set-TextValue $ws.Cells.Item(31, 4) '4.305'
$ws.Cells.Item(31, 5).Value = '  +0.70%  '
Set-TextValue $ws.Cells.Item(32, 4) '4.021'
$ws.Cells.Item(32, 5).Value = '  -1.09%  '
Set-TextValue $ws.Cells.Item(33, 4) '0.05038'
$ws.Cells.Item(33, 5).Value = '  -1.42%  '
Set-TextValue $ws.Cells.Item(34, 4) '1.219'
$ws.Cells.Item(34, 5).Value = '  +11.17%  '
Set-TextValue $ws.Cells.Item(35, 4) '0.7562'
$ws.Cells.Item(35, 5).Value = '  +5.22%  '
Set-TextValue $ws.Cells.Item(36, 4) '2.701'
$ws.Cells.Item(36, 5).Value = '  -0.31%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.01839'
$ws.Cells.Item(37, 5).Value = '  +2.26%  '
Set-TextValue $ws.Cells.Item(38, 4) '2.619'
$ws.Cells.Item(38, 5).Value = '  -0.80%  '
Set-TextValue $ws.Cells.Item(39, 4) '2.089'
$ws.Cells.Item(39, 5).Value = '  +2.25%  '
Set-TextValue $ws.Cells.Item(40, 4) '0.9053'
$ws.Cells.Item(40, 5).Value = '  -1.47%  '
Set-TextValue $ws.Cells.Item(41, 4) '5.947'
$ws.Cells.Item(41, 5).Value = '  +2.46%  '
Set-TextValue $ws.Cells.Item(42, 4) '106.85'
$ws.Cells.Item(42, 5).Value = '  +0.60%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.4294'
$ws.Cells.Item(43, 5).Value = '  +0.11%  '
$ws.Cells.Item(44, 5).Value = '  +0.41%  '
Set-TextValue $ws.Cells.Item(45, 4) '7.461'
$ws.Cells.Item(45, 5).Value = '  +0.77%  '
Set-TextValue $ws.Cells.Item(46, 4) '65.32'
$ws.Cells.Item(46, 5).Value = '  -1.53%  '
$ws.Cells.Item(47, 5).Value = '  -1.43%  '
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Cells.Item(48, 4) '1.485'
$ws.Cells.Item(48, 5).Value = '  +8.29%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(49, 4) '9.019'
$ws.Cells.Item(49, 5).Value = '  -0.64%  '
Set-TextValue $ws.Cells.Item(50, 4) '34.19'
$ws.Cells.Item(50, 5).Value = '  +0.52%  '
$ws.Cells.Item(51, 2).Value = 'Decentraland'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Cells.Item(51, 4) '0.3893'
$ws.Cells.Item(51, 5).Value = '  +1.76%  '
